$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.161.27"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.372.39"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.77"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.68"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.97%  "
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +4.13%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.22"
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.79"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.738.91"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.358.09"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.801"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.169.62"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0889"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.35"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.53"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("E28").Value = "  +5.04%  "
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.76"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  +10.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0732"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "128.70"
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.32"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.25"
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.07"
$ws.Range("E42").Value = "  -5.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.930.12"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0279"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.74"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  -8.86%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.599.19"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.51"
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "71.45"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.54"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.13"
$ws.Range("E51").Value = "  +0.41%  "
